$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing contingency ("extr") rows 8-15 with new C/D/E values ---
# Row 8 (extr1)
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# Row 9 (extr2)
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Row 10 (extr3)
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# Row 11 (extr4)
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# Row 12 (extr5)
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $true

# Row 13 (extr6)
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# Row 14 (extr7)
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

# Row 15 (extr8)
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# --- Append two new line rows (line7, line8) ---
# Row 16 (line7)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "line7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# Row 17 (line8)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "line8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true

# Copy the formatting (font/border/alignment) of the existing index column cells
# onto the two newly appended index cells so they match the rest of column A.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
